$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" date column (C) for rows 2-5 from 45221 (2023-10-22) to 45224 (2023-10-25)
$ws.Range("C2:C5").Value = 45224
